$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44308
$ws.Range("N2").Value = 280000
$ws.Range("O2").Value = 280000
$ws.Range("P2").Value = 280000
$ws.Range("R2").Value = 'Región Metropolitana'
$ws.Range("S2").Value = 800
$ws.Range("D3").Value = 44167
$ws.Range("M3").Value = 140
$ws.Range("N3").Value = 9800
$ws.Range("O3").Value = 9800
$ws.Range("P3").Value = 9800
$ws.Range("Q3").Value = '$/caja 14 kilos empedrada'
$ws.Range("R3").Value = 'Región de O''Higgins'
$ws.Range("S3").Value = 700
$ws.Range("T3").Value = 14
$ws.Range("D4").Value = 44376
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 180000
$ws.Range("O4").Value = 180000
$ws.Range("P4").Value = 180000
$ws.Range("R4").Value = 'Hijuelas'
$ws.Range("S4").Value = 514
$ws.Range("D5").Value = 44376
$ws.Range("L5").Value = 'Segunda'
$ws.Range("M5").Value = 16
$ws.Range("R5").Value = 'Provincia de Limarí'
$ws.Range("D6").Value = 44585
$ws.Range("N6").Value = 200000
$ws.Range("O6").Value = 200000
$ws.Range("P6").Value = 200000
$ws.Range("S6").Value = 571
$ws.Range("D7").Value = 44208
$ws.Range("M7").Value = 16
$ws.Range("N7").Value = 180000
$ws.Range("O7").Value = 180000
$ws.Range("P7").Value = 180000
$ws.Range("S7").Value = 514
$ws.Range("D8").Value = 44505
$ws.Range("M8").Value = 15
$ws.Range("N8").Value = 150000
$ws.Range("O8").Value = 150000
$ws.Range("P8").Value = 150000
$ws.Range("R8").Value = 'Provincia de Quillota'
$ws.Range("S8").Value = 429
$ws.Range("D9").Value = 44574
$ws.Range("M9").Value = 16
$ws.Range("N9").Value = 200000
$ws.Range("O9").Value = 200000
$ws.Range("P9").Value = 200000
$ws.Range("Q9").Value = '$/bins (350 kilos)'
$ws.Range("R9").Value = 'Región Metropolitana'
$ws.Range("S9").Value = 571
$ws.Range("T9").Value = 350
$ws.Range("D10").Value = 44657
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 160000
$ws.Range("O10").Value = 160000
$ws.Range("P10").Value = 160000
$ws.Range("R10").Value = 'Región de O''Higgins'
$ws.Range("S10").Value = 457
$ws.Range("D11").Value = 44631
$ws.Range("M11").Value = 12
$ws.Range("D12").Value = 44631
$ws.Range("L12").Value = 'Primera'
$ws.Range("M12").Value = 15
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("D13").Value = 44631
$ws.Range("L13").Value = 'Segunda'
$ws.Range("M13").Value = 15
$ws.Range("N13").Value = 170000
$ws.Range("O13").Value = 170000
$ws.Range("P13").Value = 170000
$ws.Range("R13").Value = 'Región de O''Higgins'
$ws.Range("S13").Value = 486
$ws.Range("D14").Value = 44193
$ws.Range("K14").Value = 'Start Ruby'
$ws.Range("M14").Value = 8
$ws.Range("N14").Value = 150000
$ws.Range("O14").Value = 150000
$ws.Range("P14").Value = 150000
$ws.Range("R14").Value = 'Región Metropolitana'
$ws.Range("S14").Value = 429
$ws.Range("D15").Value = 44446
$ws.Range("M15").Value = 14
$ws.Range("N15").Value = 150000
$ws.Range("O15").Value = 160000
$ws.Range("P15").Value = 155000
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 443
$ws.Range("D16").Value = 44189
$ws.Range("M16").Value = 16
$ws.Range("N16").Value = 150000
$ws.Range("O16").Value = 150000
$ws.Range("P16").Value = 150000
$ws.Range("R16").Value = 'Provincia de Limarí'
$ws.Range("S16").Value = 429
$ws.Range("D17").Value = 44627
$ws.Range("M17").Value = 6
$ws.Range("N17").Value = 240000
$ws.Range("O17").Value = 240000
$ws.Range("P17").Value = 240000
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 686
$ws.Range("D18").Value = 44196
$ws.Range("K18").Value = 'Red Blush'
$ws.Range("N18").Value = 130000
$ws.Range("O18").Value = 130000
$ws.Range("P18").Value = 130000
$ws.Range("R18").Value = 'Provincia de Limarí'
$ws.Range("S18").Value = 371
$ws.Range("D19").Value = 44400
$ws.Range("M19").Value = 140
$ws.Range("N19").Value = 9800
$ws.Range("O19").Value = 9800
$ws.Range("P19").Value = 9800
$ws.Range("Q19").Value = '$/caja 14 kilos empedrada'
$ws.Range("S19").Value = 700
$ws.Range("T19").Value = 14
$ws.Range("D20").Value = 44610
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 12
$ws.Range("N20").Value = 190000
$ws.Range("O20").Value = 190000
$ws.Range("P20").Value = 190000
$ws.Range("R20").Value = 'Región Metropolitana'
$ws.Range("S20").Value = 543
$ws.Range("D21").Value = 44356
$ws.Range("M21").Value = 24
$ws.Range("N21").Value = 200000
$ws.Range("O21").Value = 230000
$ws.Range("P21").Value = 215000
$ws.Range("R21").Value = 'Región Metropolitana'
$ws.Range("S21").Value = 614
$ws.Range("D22").Value = 44363
$ws.Range("L22").Value = 'Primera'
$ws.Range("M22").Value = 20
$ws.Range("N22").Value = 200000
$ws.Range("O22").Value = 230000
$ws.Range("P22").Value = 215000
$ws.Range("R22").Value = 'Provincia de Limarí'
$ws.Range("S22").Value = 614
$ws.Range("D23").Value = 44648
$ws.Range("M23").Value = 15
$ws.Range("N23").Value = 180000
$ws.Range("O23").Value = 180000
$ws.Range("P23").Value = 180000
$ws.Range("R23").Value = 'Región Metropolitana'
$ws.Range("S23").Value = 514
$ws.Range("D24").Value = 44645
$ws.Range("M24").Value = 24
$ws.Range("N24").Value = 170000
$ws.Range("P24").Value = 175000
$ws.Range("S24").Value = 500
$ws.Range("D25").Value = 44312
$ws.Range("L25").Value = 'Segunda'
$ws.Range("M25").Value = 10
$ws.Range("N25").Value = 330000
$ws.Range("O25").Value = 330000
$ws.Range("P25").Value = 330000
$ws.Range("S25").Value = 943
$ws.Range("D26").Value = 44586
$ws.Range("M26").Value = 20
$ws.Range("N26").Value = 180000
$ws.Range("O26").Value = 180000
$ws.Range("P26").Value = 180000
$ws.Range("Q26").Value = '$/bins (350 kilos)'
$ws.Range("R26").Value = 'Provincia de Quillota'
$ws.Range("S26").Value = 514
$ws.Range("T26").Value = 350
$ws.Range("D27").Value = 44586
$ws.Range("L27").Value = 'Segunda'
$ws.Range("M27").Value = 12
$ws.Range("N27").Value = 140000
$ws.Range("O27").Value = 140000
$ws.Range("P27").Value = 140000
$ws.Range("R27").Value = 'Provincia de Quillota'
$ws.Range("S27").Value = 400
$ws.Range("D28").Value = 44511
$ws.Range("M28").Value = 24
$ws.Range("N28").Value = 140000
$ws.Range("P28").Value = 145000
$ws.Range("S28").Value = 414
$ws.Range("D29").Value = 44195
$ws.Range("M29").Value = 20
$ws.Range("N29").Value = 200000
$ws.Range("O29").Value = 210000
$ws.Range("P29").Value = 206000
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 589
$ws.Range("D30").Value = 44201
$ws.Range("L30").Value = 'Especial'
$ws.Range("M30").Value = 8
$ws.Range("R30").Value = 'Región de O''Higgins'
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 16
$ws.Range("N31").Value = 170000
$ws.Range("O31").Value = 170000
$ws.Range("P31").Value = 170000
$ws.Range("S31").Value = 486
$ws.Range("D32").Value = 44616
$ws.Range("M32").Value = 24
$ws.Range("N32").Value = 200000
$ws.Range("O32").Value = 200000
$ws.Range("P32").Value = 200000
$ws.Range("R32").Value = 'Región Metropolitana'
$ws.Range("S32").Value = 571
$ws.Range("D33").Value = 44389
$ws.Range("L33").Value = 'Especial'
$ws.Range("M33").Value = 18
$ws.Range("N33").Value = 200000
$ws.Range("O33").Value = 200000
$ws.Range("P33").Value = 200000
$ws.Range("R33").Value = 'Provincia de Quillota'
$ws.Range("S33").Value = 571
$ws.Range("D34").Value = 44641
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 16
$ws.Range("N34").Value = 180000
$ws.Range("O34").Value = 180000
$ws.Range("P34").Value = 180000
$ws.Range("S34").Value = 514
$ws.Range("D35").Value = 44309
$ws.Range("M35").Value = 16
$ws.Range("N35").Value = 350000
$ws.Range("O35").Value = 350000
$ws.Range("P35").Value = 350000
$ws.Range("R35").Value = 'Región Metropolitana'
$ws.Range("S35").Value = 1000
$ws.Range("D36").Value = 44609
$ws.Range("L36").Value = 'Primera'
$ws.Range("M36").Value = 24
$ws.Range("N36").Value = 190000
$ws.Range("O36").Value = 190000
$ws.Range("P36").Value = 190000
$ws.Range("S36").Value = 543
